# Error Calculations and Plots
# Apply edits to the missing_data worksheet:
# 1. Remove rows 26 (RM 232) and 28 (SC 92) entirely (shifting rows up)
# 2. Fill in / clear several individual "missing data" cells in column D and E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete rows for "RM 232" (row 26) and "SC 92" (row 28) ---
# Delete from the bottom up so earlier row numbers remain valid.
$ws.Range("A28:F28").EntireRow.Delete() | Out-Null
$ws.Range("A26:F26").EntireRow.Delete() | Out-Null

# --- Individual cell edits (after the row deletions above) ---
# Row 3 (RM 8): D3 becomes missing
$ws.Range("D3").Value = $null

# Row 5 (RM 14): E5 becomes missing
$ws.Range("E5").Value = $null

# Row 8 (RM 38): E8 is filled in
$ws.Range("E8").Value = -6.6

# Row 10 (RM 52 a): E10 is filled in
$ws.Range("E10").Value = -6.1

# Row 12 (RM 81): E12 becomes missing
$ws.Range("E12").Value = $null

# Row 15 (RM 95): E15 is filled in
$ws.Range("E15").Value = -8.4

# Row 18 (RM 120): E18 becomes missing
$ws.Range("E18").Value = $null

# Row 19 (RM 125): E19 becomes missing
$ws.Range("E19").Value = $null

# Row 25 (RM 145): E25 is filled in
$ws.Range("E25").Value = -7.1

# --- Row 26 (was row 27, SC 5): B26 is filled in ---
$ws.Range("B26").Value = -20.2

# --- Row 27 (was row 29, SC 101): B27 becomes missing ---
$ws.Range("B27").Value = $null

# --- Row 29 (was row 31, SC 119): E29 becomes missing ---
$ws.Range("E29").Value = $null

# --- Row 33 (was row 35, SC 232): B33 and D33 are filled in ---
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1
